# Translate Ukrainian sheet names and shared strings to English.

$wb = $excel.ActiveWorkbook

# Map of old (Ukrainian) text -> new (English) text.
$map = @{
    "День 1" = "Day 1";
    "День 2" = "Day 2";
    "Заплив №1 - Жінки" = "Swim №1 - Females";
    "Дистанція №1 - 100 Вільний стиль" = "Distance №1 - 100 Freestyle";
    "Учасник" = "Member";
    "Вікова група" = "Age group";
    "Команда" = "Team";
    "Місто" = "City";
    "Час" = "Time";
    "Доріжка" = "Track";
    "Заплив №2 - Чоловіки" = "Swim №2 - Males";
    "Заплив №3 - Чоловіки" = "Swim №3 - Males";
    "Заплив №4 - Жінки" = "Swim №4 - Females";
    "Дистанція №2 - 50 Брас" = "Distance №2 - 50 Breaststroke";
    "Заплив №5 - Чоловіки" = "Swim №5 - Males";
    "Заплив №6 - Жінки" = "Swim №6 - Females";
    "Дистанція №3 - 100 Баттерфлай" = "Distance №3 - 100 Butterfly";
    "Заплив №7 - Чоловіки" = "Swim №7 - Males";
    "Заплив №8 - Жінки" = "Swim №8 - Females";
    "Дистанція №4 - 100 Комплексне плавання" = "Distance №4 - 100 Dolphin kick";
    "Заплив №9 - Чоловіки" = "Swim №9 - Males";
    "Заплив №3 - Жінки" = "Swim №3 - Females";
    "Дистанція №2 - 50 Вільний стиль" = "Distance №2 - 50 Freestyle";
    "Заплив №4 - Чоловіки" = "Swim №4 - Males";
}

# Replace the text in every used cell of every worksheet whenever it
# matches a key in the map (covers titles, headers, distance names).
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        $val = $cell.Value()
        if ($val -ne $null -and $map.ContainsKey([string]$val)) {
            $cell.Value = $map[[string]$val]
        }
    }
}

# Rename the worksheets themselves.
foreach ($ws in $wb.Worksheets) {
    if ($map.ContainsKey($ws.Name)) {
        $ws.Name = $map[$ws.Name]
    }
}
